# "article 86 is live"
#
# The blog-article rotation on row 7 shifts forward by one slot:
#   - article 84 is retired (its "ser: 84" card is dropped),
#   - the card that used to say "ser: 85" now lives in the slot that used to
#     hold "ser: 84" (cell E7),
#   - a brand-new "ser: 86" card takes the slot that used to hold "ser: 85"
#     (cell C7).
# Cell I7 (article 83) is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 85"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 86"

# Scroll the view up one row (topLeftCell B7 -> B6) while keeping the
# current selection (I7) untouched.
$ws.Application.ActiveWindow.ScrollRow = 6
